# Repull data, push all data, mean calculation
# Update the dSF column (column F) values for the affected rows to reflect
# the repulled / recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F5"  = -1
    "F6"  = -3
    "F10" = 4
    "F13" = 0
    "F15" = -3
    "F16" = 4
    "F18" = 0
    "F20" = -3
    "F25" = -5
    "F27" = 0
    "F28" = 0
    "F30" = -7
    "F33" = 5
    "F39" = -3
    "F41" = 0
    "F47" = -7
    "F50" = -9
    "F53" = -3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
